$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" column values read as plain numeric text (e.g. "313.41",
# "68.20", "0.000008863"). A direct $range.Value = "..." assignment lets Excel
# auto-detect these as real numbers, which silently mangles them (drops
# trailing zeros, flips to scientific notation, etc). Writing them as a
# text-returning formula and then pasting-special as values keeps the exact
# literal string without leaving a formula behind or touching cell styles.
function Set-TextValue($addr, $val) {
  $ws.Range($addr).Formula = '="' + $val + '"'
  $ws.Range($addr).Copy()
  $ws.Range($addr).PasteSpecial(-4163)
}

$ws.Range("D2").Value = '27.850.02'
$ws.Range("E2").Value = '  +2.33%  '
$ws.Range("D3").Value = '1.874.45'
$ws.Range("E3").Value = '  +0.84%  '
$ws.Range("E4").Value = '  -0.78%  '
Set-TextValue "D5" "313.41"
$ws.Range("E5").Value = '  +0.55%  '
$ws.Range("E6").Value = '  -0.71%  '
Set-TextValue "D7" "0.4827"
$ws.Range("E7").Value = '  +0.74%  '
Set-TextValue "D8" "0.3819"
$ws.Range("E8").Value = '  +2.61%  '
Set-TextValue "D9" "0.07372"
$ws.Range("E9").Value = '  +0.86%  '
Set-TextValue "D10" "0.9405"
$ws.Range("E10").Value = '  +0.40%  '
Set-TextValue "D11" "21.01"
$ws.Range("E11").Value = '  +3.91%  '
Set-TextValue "D12" "0.07787"
$ws.Range("E12").Value = '  -1.07%  '
$ws.Range("D13").Value = '1.887.44'
$ws.Range("E13").Value = '  +0.90%  '
Set-TextValue "D14" "5.517"
$ws.Range("E14").Value = '  +1.78%  '
Set-TextValue "D15" "6.615"
$ws.Range("E15").Value = '  +1.09%  '
Set-TextValue "D16" "91.36"
$ws.Range("E16").Value = '  +1.15%  '
$ws.Range("E17").Value = '  -0.85%  '
Set-TextValue "D18" "0.000008863"
$ws.Range("E18").Value = '  +1.29%  '
$ws.Range("E19").Value = '  -0.72%  '
$ws.Range("D20").Value = '27.876.06'
$ws.Range("E20").Value = '  +2.28%  '
$ws.Range("E21").Value = '  +1.10%  '
Set-TextValue "D22" "5.124"
$ws.Range("E22").Value = '  +0.54%  '
$ws.Range("D23").Value = '2.117.33'
$ws.Range("E23").Value = '  +0.41%  '
Set-TextValue "D24" "10.89"
$ws.Range("E24").Value = '  +2.17%  '
$ws.Range("E25").Value = '  -0.40%  '
Set-TextValue "D26" "157.63"
$ws.Range("E26").Value = '  +2.54%  '
Set-TextValue "D27" "18.56"
$ws.Range("E27").Value = '  +0.29%  '
Set-TextValue "D28" "2.042"
$ws.Range("E28").Value = '  +2.36%  '
Set-TextValue "D29" "116.06"
$ws.Range("E29").Value = '  +0.30%  '
Set-TextValue "D30" "4.977"
$ws.Range("E30").Value = '  +0.68%  '
Set-TextValue "D31" "0.08888"
$ws.Range("E31").Value = '  -0.03%  '
Set-TextValue "D32" "3.341"
$ws.Range("E32").Value = '  -0.14%  '
Set-TextValue "D33" "1.223"
$ws.Range("E33").Value = '  +3.33%  '
Set-TextValue "D34" "0.7724"
$ws.Range("E34").Value = '  +4.27%  '
Set-TextValue "D35" "4.651"
$ws.Range("E35").Value = '  +1.36%  '
Set-TextValue "D36" "2.738"
$ws.Range("E36").Value = '  +2.19%  '
Set-TextValue "D37" "1.126"
$ws.Range("E37").Value = '  +0.28%  '
Set-TextValue "D38" "0.02047"
$ws.Range("E38").Value = '  +1.37%  '
Set-TextValue "D39" "0.5629"
$ws.Range("E39").Value = '  +5.41%  '
Set-TextValue "D40" "0.05374"
$ws.Range("E40").Value = '  +2.26%  '
$ws.Range("E41").Value = '  +0.13%  '
Set-TextValue "D42" "7.059"
$ws.Range("E42").Value = '  -0.61%  '
Set-TextValue "D43" "8.543"
$ws.Range("E43").Value = '  +2.49%  '
Set-TextValue "D44" "0.1530"
$ws.Range("E44").Value = '  +0.17%  '
Set-TextValue "D45" "0.4884"
$ws.Range("E45").Value = '  +2.14%  '
Set-TextValue "D46" "10.69"
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D47" "105.64"
$ws.Range("E47").Value = '  +2.87%  '
$ws.Range("B48").Value = 'PaxDollar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue "D48" "1.012"
$ws.Range("E48").Value = '  -0.69%  '
Set-TextValue "D49" "1.666"
$ws.Range("E49").Value = '  +2.01%  '
Set-TextValue "D50" "68.20"
Set-TextValue "D51" "0.06122"
$ws.Range("E51").Value = '  +0.72%  '
